# Apply the "remove RG, keep only CPF" wording change and the matching
# shrink of the certificate text box (CertificateNR15 model), as
# described by the source diff.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# The paragraph lives in the "Rectangle 5" shape (2nd shape on the slide).
$shape = $s.Shapes.Item(2)
$tf = $shape.TextFrame
$tr = $tf.TextRange

# --- Text edits -----------------------------------------------------
# Work from the end of the text range towards the start so that the
# character offsets used below stay valid after each edit shrinks the
# text.

# 1) Last paragraph: "São Carlos" + ", " + "{{DATA}}" -> merge the first
#    two runs into a single "São Carlos, " run (text unchanged, only the
#    run split disappears).
$tr.Characters(194, 12).Text = "São Carlos, "

# 2) First paragraph: "...portador do RG nº {{RG}} e CPF nº {{CPF}}..."
#    becomes "...portador do CPF nº {{CPF}}...". Edit the three affected
#    runs in place, right to left, so each keeps referring to valid,
#    still-unshifted offsets.
$tr.Characters(53, 9).Text = "nº "          # was "e CPF nº "
$tr.Characters(46, 7).Text = "do CPF "      # was "{{RG}} "
$tr.Characters(28, 18).Text = "portador "   # was "portador do RG nº "

# --- Resize the (auto-fit) text box ----------------------------------
# Less text now fits in a shorter box; height shrinks accordingly while
# position/width stay the same (936876, 2366189, 7571874 EMU untouched).
$shape.Height = 155.1
